$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to English short column names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Normalize municipality/state names: title-case Spanish connector words
# (de/del/la/las/los/el/y) to match the new naming convention, and fix
# the "MonteMorelos" typo to "Montemorelos"
$ws.Range("B4").Value = 'Pabellón De Arteaga'
$ws.Range("B21").Value = 'Chiapa De Corzo'
$ws.Range("B24").Value = 'Comitán De Domínguez'
$ws.Range("B39").Value = 'Mazapa De Madero'
$ws.Range("B42").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B49").Value = 'Salto De Agua'
$ws.Range("B50").Value = 'San Cristóbal De Las Casas'
$ws.Range("B78").Value = 'Hidalgo Del Parral'
$ws.Range("B83").Value = 'San Francisco De Borja'
$ws.Range("A105").Value = 'Ciudad De México'
$ws.Range("B126").Value = 'Pánuco De Coronado'
$ws.Range("A131").Value = 'Estado De México'
$ws.Range("B131").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B132").Value = 'Almoloya De Alquisiras'
$ws.Range("B140").Value = 'Ecatepec De Morelos'
$ws.Range("B144").Value = 'Ixtapan De La Sal'
$ws.Range("B148").Value = 'Naucalpan De Juárez'
$ws.Range("B151").Value = 'San Felipe Del Progreso'
$ws.Range("B159").Value = 'Tlalnepantla De Baz'
$ws.Range("B163").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B169").Value = 'Apaseo El Alto'
$ws.Range("B170").Value = 'Apaseo El Grande'
$ws.Range("B176").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B183").Value = 'Purísima Del Rincón'
$ws.Range("B186").Value = 'San Diego De La Unión'
$ws.Range("B189").Value = 'San Luis De La Paz'
$ws.Range("B191").Value = 'Valle De Santiago'
$ws.Range("B196").Value = 'Acapulco De Juárez'
$ws.Range("B197").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B199").Value = 'Atoyac De Álvarez'
$ws.Range("B200").Value = 'Ayutla De Los Libres'
$ws.Range("B202").Value = 'Chilapa De Álvarez'
$ws.Range("B203").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B204").Value = 'Coyuca De Catalán'
$ws.Range("B205").Value = 'Cuetzala Del Progreso'
$ws.Range("B207").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B209").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B212").Value = 'Mártir De Cuilapan'
$ws.Range("B218").Value = 'Taxco De Alarcón'
$ws.Range("B220").Value = 'Técpan De Galeana'
$ws.Range("B222").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B226").Value = 'Tlapa De Comonfort'
$ws.Range("B233").Value = 'Agua Blanca De Iturbide'
$ws.Range("B243").Value = 'Huejutla De Reyes'
$ws.Range("B248").Value = 'Pachuca De Soto'
$ws.Range("B251").Value = 'Progreso De Obregón'
$ws.Range("B254").Value = 'Santiago De Anaya'
$ws.Range("B255").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B258").Value = 'Tenango De Doria'
$ws.Range("B260").Value = 'Tepehuacán De Guerrero'
$ws.Range("B261").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B262").Value = 'Tezontepec De Aldama'
$ws.Range("B266").Value = 'Tula De Allende'
$ws.Range("B267").Value = 'Tulancingo De Bravo'
$ws.Range("B269").Value = 'Zacualtipán De Ángeles'
$ws.Range("B274").Value = 'Atotonilco El Alto'
$ws.Range("B281").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B283").Value = 'Jilotlán De Los Dolores'
$ws.Range("B285").Value = 'Lagos De Moreno'
$ws.Range("B288").Value = 'Ojuelos De Jalisco'
$ws.Range("B290").Value = 'San Juan De Los Lagos'
$ws.Range("B291").Value = 'San Miguel El Alto'
$ws.Range("B293").Value = 'Tamazula De Gordiano'
$ws.Range("B295").Value = 'Tepatitlán De Morelos'
$ws.Range("B297").Value = 'Tizapán El Alto'
$ws.Range("B300").Value = 'Yahualica De González Gallo'
$ws.Range("B302").Value = 'Zapotlán El Grande'
$ws.Range("B358").Value = 'Montemorelos'
$ws.Range("B360").Value = 'San Nicolás De Los Garza'
$ws.Range("B363").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B365").Value = 'Capulálpam De Méndez'
$ws.Range("B367").Value = 'Coicoyán De Las Flores'
$ws.Range("B370").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B371").Value = 'Ixtlán De Juárez'
$ws.Range("B372").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B374").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B375").Value = 'Putla Villa De Guerrero'
$ws.Range("B389").Value = 'San Miguel El Grande'
$ws.Range("B399").Value = 'Santa Ana Del Valle'
$ws.Range("B407").Value = 'Santa Lucía Del Camino'
$ws.Range("B412").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B426").Value = 'Santo Domingo De Morelos'
$ws.Range("B432").Value = 'Tanetze De Zaragoza'
$ws.Range("B433").Value = 'Tataltepec De Valdés'
$ws.Range("B434").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B435").Value = 'Tlacolula De Matamoros'
$ws.Range("B436").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B437").Value = 'Zimatlán De Álvarez'
$ws.Range("B449").Value = 'Los Reyes De Juárez'
$ws.Range("B451").Value = 'Palmar De Bravo'
$ws.Range("B455").Value = 'San Salvador El Seco'
$ws.Range("B458").Value = 'Tetela De Ocampo'
$ws.Range("B471").Value = 'Amealco De Bonfil'
$ws.Range("B473").Value = 'Cadereyta De Montes'
$ws.Range("B477").Value = 'Jalpan De Serra'
$ws.Range("B479").Value = 'Pinal De Amoles'
$ws.Range("B482").Value = 'San Juan Del Río'
$ws.Range("B492").Value = 'Cerro De San Pedro'
$ws.Range("B493").Value = 'Ciudad Del Maíz'
$ws.Range("B499").Value = 'Mexquitic De Carmona'
$ws.Range("B504").Value = 'San Ciro De Acosta'
$ws.Range("B508").Value = 'Santa María Del Río'
$ws.Range("B516").Value = 'Tanquián De Escobedo'
$ws.Range("B518").Value = 'Villa De Guadalupe'
$ws.Range("B519").Value = 'Villa De La Paz'
$ws.Range("B520").Value = 'Villa De Ramos'
$ws.Range("B521").Value = 'Villa De Reyes'
$ws.Range("B536").Value = 'Nacozari De García'
$ws.Range("B563").Value = 'Soto La Marina'
$ws.Range("B572").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B577").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B580").Value = 'Amatlán De Los Reyes'
$ws.Range("B590").Value = 'Cosamaloapan De Carpio'
$ws.Range("B591").Value = 'Cosautlán De Carvajal'
$ws.Range("B599").Value = 'Hueyapan De Ocampo'
$ws.Range("B600").Value = 'Ignacio De La Llave'
$ws.Range("B611").Value = 'Lerdo De Tejada'
$ws.Range("B612").Value = 'Martínez De La Torre'
$ws.Range("B614").Value = 'Medellín De Bravo'
$ws.Range("B624").Value = 'Paso De Ovejas'
$ws.Range("B628").Value = 'Poza Rica De Hidalgo'
$ws.Range("B634").Value = 'Sayula De Alemán'
$ws.Range("B635").Value = 'Soledad De Doblado'
$ws.Range("B636").Value = 'Tatahuicapan De Juárez'
$ws.Range("B655").Value = 'Vega De Alatorre'
$ws.Range("B665").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B674").Value = 'Nochistlán De Mejía'
$ws.Range("B681").Value = 'Teúl De González Ortega'

# Correct rounding of the percentage value in D589
$ws.Range("D589").Value = 0.009138381201044389

# Remove the trailing footnote/metadata rows (689-693) that are not part of
# the tabular data
$ws.Range("A689:A693").EntireRow.Delete()
